{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph introducing the VR game criteria list\n// (\"Cr\u00e9er un jeu vid\u00e9o en r\u00e9alit\u00e9 virtuelle ...\").\nlet introIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Cr\u00e9er un jeu vid\u00e9o\") === 0) {\n    introIndex = i;\n    break;\n  }\n}\n\nif (introIndex !== -1) {\n  // Remove the blank paragraph plus the three bullet paragraphs that used\n  // to follow the intro (the \"objectif clair\", \"contr\u00f4les\" and\n  // \"L'objectif du jeu est clair (10%)\" items), leaving the remaining\n  // bullets (\"Utilisation judicieuse des menus\", \"Clart\u00e9 et documentation\n  // du code\") directly after the intro paragraph.\n  const toRemove = [\n    paragraphs.items[introIndex + 1],\n    paragraphs.items[introIndex + 2],\n    paragraphs.items[introIndex + 3],\n    paragraphs.items[introIndex + 4],\n  ];\n\n  // Delete from the last one to the first one so earlier deletions don't\n  // invalidate the references of the ones still pending.\n  for (let i = toRemove.length - 1; i >= 0; i--) {\n    toRemove[i].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that introduces the criteria list (\"Cr\u00e9er un jeu\n# vid\u00e9o en r\u00e9alit\u00e9 virtuelle ...\"), then remove the blank paragraph and the\n# three bullet paragraphs that used to follow it (the \"objectif clair\",\n# \"contr\u00f4les\" and \"L'objectif du jeu est clair (10%)\" items), leaving the\n# remaining bullets (\"Utilisation judicieuse des menus\", \"Clart\u00e9 et\n# documentation du code\") intact right after the intro paragraph.\n\n$introText = \"Cr\" + [char]0x00E9 + \"er un jeu vid\" + [char]0x00E9 + \"o\"\n$startPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"$introText*\") {\n        $startPara = $p\n        break\n    }\n}\n\n$firstToRemove = $startPara.Next()\n$lastToRemove = $firstToRemove.Next().Next().Next()\n\n$range = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)\n$range.Delete()\n"}
